# Edit script: update NATMI LR-pair values ("Natmi following Dr Hou advice")
# Updates columns E,G,H,I,J,K,M,N,O,P,Q,R,S,T for rows 2-16 on the active sheet,
# reflecting recomputed ligand/receptor-expressing cell counts (1 -> 3) and the
# derived expression / specificity statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2  = @{ "E"="3"; "G"="0.881601"; "H"="2.644803"; "I"="0.02072192623875441"; "J"="0.02072192623875442"; "K"="3"; "M"="22.34478233333333"; "N"="67.034347"; "O"="0.03165884810812076"; "P"="0.03165884810812076"; "Q"="19.699182449849"; "R"="177.292642048641"; "S"="0.000656032315300408"; "T"="0.0006560323153004081" }
  3  = @{ "E"="3"; "G"="0.881601"; "H"="2.644803"; "I"="0.02072192623875441"; "J"="0.02072192623875442"; "K"="3"; "M"="92.44713066666667"; "N"="277.341392"; "O"="0.1309822411400946"; "P"="0.1309822411400946"; "Q"="81.50148284286399"; "R"="733.5133455857759"; "S"="0.002714204339491785"; "T"="0.002714204339491785" }
  4  = @{ "E"="3"; "G"="0.881601"; "H"="2.644803"; "I"="0.02072192623875441"; "J"="0.02072192623875442"; "K"="3"; "M"="243.96462"; "N"="731.89386"; "O"="0.3456573768818275"; "P"="0.3456573768818275"; "Q"="215.07945295662"; "R"="1935.71507660958"; "S"="0.007162686667626564"; "T"="0.007162686667626565" }
  5  = @{ "E"="3"; "G"="0.881601"; "H"="2.644803"; "I"="0.02072192623875441"; "J"="0.02072192623875442"; "K"="3"; "M"="281.5837096666667"; "N"="844.751129"; "O"="0.3989573834764815"; "P"="0.3989573834764815"; "Q"="248.244480025843"; "R"="2234.200320232587"; "S"="0.008267165472806108"; "T"="0.008267165472806108" }
  6  = @{ "E"="3"; "G"="0.881601"; "H"="2.644803"; "I"="0.02072192623875441"; "J"="0.02072192623875442"; "K"="3"; "M"="65.45872566666667"; "N"="196.376177"; "O"="0.09274415039347572"; "P"="0.09274415039347571"; "Q"="57.708478006459"; "R"="519.376302058131"; "S"="0.00192183744352955"; "T"="0.00192183744352955" }
  7  = @{ "E"="3"; "G"="37.27504099999999"; "H"="111.825123"; "I"="0.8761453879346173"; "J"="0.8761453879346174"; "K"="3"; "M"="22.34478233333333"; "N"="67.034347"; "O"="0.03165884810812076"; "P"="0.03165884810812076"; "Q"="832.9026776110754"; "R"="7496.12409849968"; "S"="0.02773775375725259"; "T"="0.02773775375725259" }
  8  = @{ "E"="3"; "G"="37.27504099999999"; "H"="111.825123"; "I"="0.8761453879346173"; "J"="0.8761453879346174"; "K"="3"; "M"="92.44713066666667"; "N"="277.341392"; "O"="0.1309822411400946"; "P"="0.1309822411400946"; "Q"="3445.970585932357"; "R"="31013.73527339121"; "S"="0.1147594864762338"; "T"="0.1147594864762338" }
  9  = @{ "E"="3"; "G"="37.27504099999999"; "H"="111.825123"; "I"="0.8761453879346173"; "J"="0.8761453879346174"; "K"="3"; "M"="243.96462"; "N"="731.89386"; "O"="0.3456573768818275"; "P"="0.3456573768818275"; "Q"="9093.791213049419"; "R"="81844.12091744477"; "S"="0.3028461165605909"; "T"="0.302846116560591" }
  10 = @{ "E"="3"; "G"="37.27504099999999"; "H"="111.825123"; "I"="0.8761453879346173"; "J"="0.8761453879346174"; "K"="3"; "M"="281.5837096666667"; "N"="844.751129"; "O"="0.3989573834764815"; "P"="0.3989573834764815"; "Q"="10496.04432275709"; "R"="94464.39890481386"; "S"="0.3495446715153818"; "T"="0.3495446715153818" }
  11 = @{ "E"="3"; "G"="37.27504099999999"; "H"="111.825123"; "I"="0.8761453879346173"; "J"="0.8761453879346174"; "K"="3"; "M"="65.45872566666667"; "N"="196.376177"; "O"="0.09274415039347572"; "P"="0.09274415039347571"; "Q"="2439.976683032752"; "R"="21959.79014729477"; "S"="0.08125735962515827"; "T"="0.08125735962515827" }
  12 = @{ "E"="3"; "G"="4.387713666666667"; "H"="13.163141"; "I"="0.1031326858266283"; "J"="0.1031326858266283"; "K"="3"; "M"="22.34478233333333"; "N"="67.034347"; "O"="0.03165884810812076"; "P"="0.03165884810812076"; "Q"="98.04250682265854"; "R"="882.382561403927"; "S"="0.003265062035567764"; "T"="0.003265062035567764" }
  13 = @{ "E"="3"; "G"="4.387713666666667"; "H"="13.163141"; "I"="0.1031326858266283"; "J"="0.1031326858266283"; "K"="3"; "M"="92.44713066666667"; "N"="277.341392"; "O"="0.1309822411400946"; "P"="0.1309822411400946"; "Q"="405.6315386702524"; "R"="3650.683848032272"; "S"="0.01350855032436905"; "T"="0.01350855032436905" }
  14 = @{ "E"="3"; "G"="4.387713666666667"; "H"="13.163141"; "I"="0.1031326858266283"; "J"="0.1031326858266283"; "K"="3"; "M"="243.96462"; "N"="731.89386"; "O"="0.3456573768818275"; "P"="0.3456573768818275"; "Q"="1070.44689735714"; "R"="9634.02207621426"; "S"="0.03564857365360997"; "T"="0.03564857365360997" }
  15 = @{ "E"="3"; "G"="4.387713666666667"; "H"="13.163141"; "I"="0.1031326858266283"; "J"="0.1031326858266283"; "K"="3"; "M"="281.5837096666667"; "N"="844.751129"; "O"="0.3989573834764815"; "P"="0.3989573834764815"; "Q"="1235.508691215132"; "R"="11119.57822093619"; "S"="0.04114554648829364"; "T"="0.04114554648829363" }
  16 = @{ "E"="3"; "G"="4.387713666666667"; "H"="13.163141"; "I"="0.1031326858266283"; "J"="0.1031326858266283"; "K"="3"; "M"="65.45872566666667"; "N"="196.376177"; "O"="0.09274415039347572"; "P"="0.09274415039347571"; "Q"="287.2141452102175"; "R"="2584.927306891957"; "S"="0.009564953324787897"; "T"="0.009564953324787897" }
}

foreach ($rowNum in $data.Keys) {
    $rowVals = $data[$rowNum]
    foreach ($col in $rowVals.Keys) {
        $addr = "$col$rowNum"
        $ws.Range($addr).Value = [double]$rowVals[$col]
    }
}
